$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Used By Process" (column D) text for rows whose content changed ---
$ws.Range("D3").Value = "SCR_-_SOP_-_Buy_a_car,5Y_Department_Budget_and_HC_Plan,EU_Reimbursement_Policy,User_wants_to_book_testdrive_online,Get_vehicles_from_CN_or_EU_hub_to_RDC,Receive_subscription_order_from_user,Send_vehicle_preparation_order_to_RDC_and_the_registration_order_to_DAD,Vehicle_manager_check_with_RDC_about_vehicle_receive_and_PDI_preparation_related_issues,Deliver_subscription_cars,Subscription_cars_delivery-related,Fellow_complete_testdrive_process,Fellow_manually_books_testdrive_for_user,Fellow_handles_a_testdrive_cancelation_or_change_request,Fellow_follows_up_testdrive_leads,Fellow_fills_damage_form_together_with_users,test_process,Fellow_test_drive_follow_up_sop,MSRP_Pricing_SOP,User Behavior Tracking System (UBTS),Market_Intelligence_Europe,Sales_Planning,Sales_Planning - Supply Planning,Sales_Planning - Sales_Planning,Sales_Planning - Fleet Planning,Sales_Planning - Supply Planning - offering structure,Sales_Planning - Supply Planning - planning to production,Sales_Planning - Supply Planning - Allocation & Re-allocation,Sales_Planning - Sales_Planning - Sales Performance Rolling Plan,Sales_Planning - Sales_Planning - Roadmap to Target(Target Streeting),Sales_Planning - Fleet Planning - stock level monitor & stock age forecast model,Sales_Planning - Fleet Planning - subscription fleet plan & monitor,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - in-fleet plan,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - de-fleet plan"
$ws.Range("D4").Value = "SCR_-_SOP_-_Activate_Home_Charger,Colleague_Request_New_System_Feature,1Y_Department_Budget_and_HC_Plan,EU_Data_Security_Control_Policy,NIO_Life_-_1900-190X - Internal Application,NIO_Life_-_APP_sales - Product launch,NIO_Life_-_APP_sales - After sales"
$ws.Range("D7").Value = "SCR_-_SOP_-_Handles_Low_Star,5_Year_Sales_Volume_Calculation,Process_Survey_Evaluation,IT_Hardware_give_out_management,Complete_Purchasing_process,EB_quartlery_5_year_budget_planning_process,Power_Procurement,Power_EU_Supply_Chain_Order_to_Delivery,Power_Procurement - Standard Process,Power_Procurement - Authorized Process,Power EU PUS Leads to Operation (Q-300-B-1 EU PUS LTO ),NIO_Life_-_Export - Export for internal applicatioin,MSRP_Pricing_SOP"
$ws.Range("D8").Value = "SCR_Complaint_Handling,5_Year_Strategy_Planning_Presentation,Monthly_Process_Review,Sickness_Leave,Maintenance of Register of Processing Activities (RoPA),5yr_planning_flow,Create_a_new_promotion_,MSRP_Pricing_SOP,User Behavior Tracking System (UBTS),Market_Intelligence_Europe,SUD_Experience_Confirmation"
$ws.Range("D10").Value = "SCR_-_SOP_-_Test_Drive_Request,User_Request_Leasing,Agent_Performance_Review,Power_EU_PUS_Leads_to_Operation_-_Power_EU_PUS_Step_by_Step"
$ws.Range("D21").Value = "Vehicle_manager_check_with_RDC_about_vehicle_receive_and_PDI_preparation_related_issues,Get_vehicles_from_CN_or_EU_hub_to_RDC,Receive_subscription_order_from_user,Send_vehicle_preparation_order_to_RDC_and_the_registration_order_to_DAD,Deliver_subscription_cars,Subscription_cars_delivery-related,User Behavior Tracking System (UBTS),Sales_Planning,Sales_Planning - Supply Planning,Sales_Planning - Sales_Planning,Sales_Planning - Fleet Planning,Sales_Planning - Supply Planning - offering structure,Sales_Planning - Supply Planning - planning to production,Sales_Planning - Supply Planning - Allocation & Re-allocation,Sales_Planning - Sales_Planning - Sales Performance Rolling Plan,Sales_Planning - Sales_Planning - Roadmap to Target(Target Streeting),Sales_Planning - Fleet Planning - stock level monitor & stock age forecast model,Sales_Planning - Fleet Planning - subscription fleet plan & monitor,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - in-fleet plan,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - de-fleet plan"
$ws.Range("D24").Value = "Fellow_handles_a_Testdrive_Process,Privacy Impact Assessment (PIA),hr test_process,ud_publishes_a_pgc,Network_development_manager_selects_design_vendor,Power_EU_Supply_Chain_Planning_&_Forecast,Power_EU_PUS_Leads_to_Operation_-_EU_Annual_Planning_and_Budget_Approval,Power_EU_PUS_Leads_to_Operation_-_Country_Annual_Planning_and_Budget_Approval,Power_EU_PUS_Leads_to_Operation_-_LTO_Process,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Leads to Contract,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation, Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Leads to Contract - Leads Collection,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Leads to Contract - Site inspection,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Leads to Contract - Quotation,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Leads to Contract - Project approval,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Leads to Contract - Contracting,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Announcement,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Announcement - Kickoff,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Announcement - Partner Nomination (if applicable),Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Construction Permit,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Grid Connection & Upgrade (if applicable),Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Civil Construction,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Installation & Commissioning,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Installation & Commissioning - Device installation,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Installation & Commissioning - Installation & Commissioning (owned by NIO Power),Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - PUS Ready for Service,Power_EU_PUS_Leads_to_Operation_-_LTO_Process - Contract to Operation - Launch,Power_EU_PUS_Leads_to_Operation_-_PUS_Operation_Approval_Process_,NIO_Life_-_Export,NIO_Life_-_Export - Planned export,NIO_Life_-_1900-190X,NIO_Life_-_1900-190X - New store opening,NIO_Life_-_1900-190X - Replenishment,NIO_Life_-_APP_sales,EU_Legal_-_Data_Protection,Europe_Market_Planning,Europe_Market_Planning - 5 yrs planning,Europe_Market_Planning - 5 yrs planning - 5 yrs sales volume planning,Europe_Market_Planning - 5 yrs planning - market expansion 5 yrs assumptions,Europe_Market_Planning - 5 yrs planning - market expansion 5 yrs assumptions - future portfolio mapping,Europe_Market_Planning - 5 yrs planning - market expansion 5 yrs assumptions - business model,Europe_Market_Planning - 5 yrs planning - market expansion 5 yrs assumptions - market entry cadence & infrasturcutre,Europe_Market_Planning - Wave2 / Market Entry PMO,Europe_Market_Planning - Governance PMO,Europe_Market_Planning - Governance PMO - Entity setup,EU_Legal_-_Data_Protection - Data Protection Risk Management"

# --- Append new system rows 25-32 ---
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "NIO_App"
$ws.Range("C25").Value = "Click Me"
$ws.Range("D25").Value = "Power_EU_Supply_Chain_Order_to_Delivery - User Orders via NIO_App"

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "VLM"
$ws.Range("C26").Value = "Click Me"
$ws.Range("D26").Value = "Sales_Planning,Sales_Planning - Sales_Planning,Sales_Planning - Fleet Planning,Sales_Planning - Supply Planning - offering structure,Sales_Planning - Supply Planning - planning to production,Sales_Planning - Supply Planning - Allocation & Re-allocation,Sales_Planning - Fleet Planning - stock level monitor & stock age forecast model,Sales_Planning - Fleet Planning - subscription fleet plan & monitor,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - in-fleet plan,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - de-fleet plan"

$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "LOPA"
$ws.Range("C27").Value = "Click Me"
$ws.Range("D27").Value = "Sales_Planning,Sales_Planning - Supply Planning,Sales_Planning - Supply Planning - offering structure,Sales_Planning - Supply Planning - planning to production,Sales_Planning - Supply Planning - Allocation & Re-allocation,Sales_Planning - Sales_Planning - Sales Performance Rolling Plan,Sales_Planning - Sales_Planning - Roadmap to Target(Target Streeting)"

$ws.Range("A28").Value = 26
$ws.Range("C28").Value = "Click Me"

$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "Excel"
$ws.Range("C29").Value = "Click Me"
$ws.Range("D29").Value = "Sales_Planning,Sales_Planning - Sales_Planning,Sales_Planning - Fleet Planning,Sales_Planning - Supply Planning - offering structure,Sales_Planning - Supply Planning - planning to production,Sales_Planning - Supply Planning - Allocation & Re-allocation,Sales_Planning - Fleet Planning - stock level monitor & stock age forecast model,Sales_Planning - Fleet Planning - subscription fleet plan & monitor,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - in-fleet plan,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - de-fleet plan"

$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "Tableau"
$ws.Range("C30").Value = "Click Me"
$ws.Range("D30").Value = "Sales_Planning,Sales_Planning - Sales_Planning,Sales_Planning - Fleet Planning,Sales_Planning - Supply Planning - offering structure,Sales_Planning - Supply Planning - planning to production,Sales_Planning - Supply Planning - Allocation & Re-allocation,Sales_Planning - Fleet Planning - stock level monitor & stock age forecast model,Sales_Planning - Fleet Planning - subscription fleet plan & monitor,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - in-fleet plan,Sales_Planning - Fleet Planning - subscription fleet plan & monitor - de-fleet plan"

$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "UCLS"
$ws.Range("C31").Value = "Click Me"
$ws.Range("D31").Value = "User_Closed-Loop"

$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "IPD++"
$ws.Range("C32").Value = "Click Me"
$ws.Range("D32").Value = "User_Closed-Loop"

# --- Apply the same formatting as the existing numbered rows to column A of the new rows ---
$newA = $ws.Range("A25:A32")
$newA.Font.Bold = $true
$newA.HorizontalAlignment = -4108
$newA.VerticalAlignment = -4160
$newA.Borders.LineStyle = 1

Write-Output "done"
